# Updated symbol list on Sat Dec 24 21:48:09 UTC 2022 with GitHub Actions
#
# This script applies the price/volume refresh captured in the diff:
#  - Column D ("Price") text updates for most rows
#  - Rows 41-43 ("Coin" B, "Link" C, "Price" D, "Volume(1h)" E) get
#    reshuffled because the coin ranking order changed (KickToken moved
#    up to row 41, BKEXToken moved to row 42, CEJI moved to row 43)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these cells hold plain text (e.g. "244.50", "0.002319") rather
# than numbers, so every new value is written with a leading apostrophe
# to force Excel to keep it as text instead of auto-converting it to a
# Number (which would drop meaningful trailing zeros / use exponent
# notation for very small values).
function Set-TextValue($address, $value) {
    $ws.Range($address).Value = "'" + $value
}

# Column D ("Price") updates that do not involve any row reordering
$priceUpdates = [ordered]@{
    "D2"  = "244.58"
    "D3"  = "21.89"
    "D4"  = "5.409"
    "D5"  = "0.06039"
    "D6"  = "3.396"
    "D7"  = "0.8144"
    "D8"  = "0.9205"
    "D9"  = "0.1439"
    "D10" = "0.07436"
    "D11" = "0.03408"
    "D12" = "0.03050"
    "D13" = "0.09426"
    "D14" = "4.009"
    "D15" = "0.001591"
    "D16" = "0.04819"
    "D17" = "0.0005946"
    "D19" = "0.004163"
    "D20" = "0.0009885"
    "D22" = "6.427"
    "D25" = "0.1324"
    "D26" = "0.00008405"
    "D27" = "0.0002903"
    "D40" = "0.03996"
    "D44" = "0.005796"
    "D45" = "0.00005246"
    "D46" = "0.00000000751"
    "D47" = "1.001"
    "D48" = "0.002320"
    "D49" = "0.00002102"
    "D50" = "0.01011"
}

foreach ($address in $priceUpdates.Keys) {
    Set-TextValue $address $priceUpdates[$address]
}

# Rows 41-43 got re-ranked: KickToken now leads (row 41), followed by
# BKEXToken (row 42) and CEJI (row 43), each with refreshed prices.
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006423"
Set-TextValue "E41" "40KickTokenKICK"

Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1077"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002779"
Set-TextValue "E43" "42CEJICEJI"
